$d = $word.ActiveDocument

# The final paragraph in the document currently holds nothing but the
# "_GoBack" bookmark. We need to:
#   1) insert a new paragraph "Problem 1:  A Cat, Parrot and a Bag of Seed"
#   2) insert a blank paragraph after it
#   3) insert a run of text ("Defining the problem: ...") into the existing
#      bookmark paragraph, ahead of the bookmark start/end tags.
# all three new paragraphs/runs go right before that bookmark paragraph.

$bookmarkPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$bookmarkRange = $bookmarkPara.Range

# Put the "Defining the problem" text straight into the bookmark paragraph,
# ahead of the bookmark, by collapsing the range to its very start first.
$bookmarkRange.Collapse(1)
$bookmarkRange.InsertBefore("Defining the problem: We have a man, a cat, a parrot and a bag of seed that need to be transported across the river.")

# Now split off two new paragraphs before that same (now-merged) paragraph:
# one blank, and one for the "Problem 1" title.
$splitRange = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$splitRange.Collapse(1)
$splitRange.InsertParagraphBefore()
$splitRange.InsertParagraphBefore()

$titleIndex = $d.Paragraphs.Count - 2
$d.Paragraphs.Item($titleIndex).Range.Text = "Problem 1:  A Cat, Parrot and a Bag of Seed"
